$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 4 and 5 correspond to 4434c678... and 91551ae0... files
foreach ($r in 4,5) {
    $overview.Range("E$r").Value = "Ready for handoff"
    $overview.Range("F$r").Value = "Ready for handoff"
    $overview.Range("G$r").Value = "2016-11-09 01:34:40"
}

# zh-cn sheet: rows 4 and 5
$zhcn.Range("C4").Value = "Ready for handoff"
$zhcn.Range("H4").Value = "2016-11-09 01:34:25"
$zhcn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2243835334a890c6cfab22c9f6022cacca9b453b/e2e/4434c678-1261-495d-8d7a-a5dda76c6919.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54aacec4a2bc01cbf8ac34e783ed6bb72c887f9b/e2e/4434c678-1261-495d-8d7a-a5dda76c6919.md."

$zhcn.Range("C5").Value = "Ready for handoff"
$zhcn.Range("H5").Value = "2016-11-09 01:34:25"
$zhcn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2243835334a890c6cfab22c9f6022cacca9b453b/e2e/91551ae0-ea99-4516-abd5-ad24bc056486.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54aacec4a2bc01cbf8ac34e783ed6bb72c887f9b/e2e/91551ae0-ea99-4516-abd5-ad24bc056486.md."

# de-de sheet: rows 4 and 5
$dede.Range("C4").Value = "Ready for handoff"
$dede.Range("H4").Value = "2016-11-09 01:34:40"
$dede.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2243835334a890c6cfab22c9f6022cacca9b453b/e2e/4434c678-1261-495d-8d7a-a5dda76c6919.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54aacec4a2bc01cbf8ac34e783ed6bb72c887f9b/e2e/4434c678-1261-495d-8d7a-a5dda76c6919.md."

$dede.Range("C5").Value = "Ready for handoff"
$dede.Range("H5").Value = "2016-11-09 01:34:40"
$dede.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2243835334a890c6cfab22c9f6022cacca9b453b/e2e/91551ae0-ea99-4516-abd5-ad24bc056486.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/54aacec4a2bc01cbf8ac34e783ed6bb72c887f9b/e2e/91551ae0-ea99-4516-abd5-ad24bc056486.md."

# Widen column P (Error Detail) on zh-cn and de-de sheets to XML width=40.
# This engine's ColumnWidth -> stored xml "width" conversion adds a fixed
# 5/6 offset, so request 39 + 1/6 to land exactly on 40.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
